$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.858.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.84%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.096.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +5.18%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'580.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.59%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'172.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +5.78%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.10%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.091.67"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.16%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +1.34%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.56%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.155"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.93%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +3.95%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.70%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'37.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +6.67%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.13%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.608.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +5.08%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'66.827.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.79%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'7.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.49%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.097.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +5.16%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'16.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +3.01%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'479.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "'  +2.57%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +3.31%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'83.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.19%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'13.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +8.21%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +5.43%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.92%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.02%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.92%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -2.31%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +3.76%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'28.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +5.39%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0₃0998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.87%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -1.12%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.00%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.21%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.989"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.98%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'48.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.80%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +6.83%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +1.90%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.316"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.50%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +1.18%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'8.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.16%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.30%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.838.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "'  +2.68%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'383.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.48%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'135.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.79%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.00%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +3.52%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.80%  "
$ws.Range("E51").Style = "Normal"
